# GSC export update: append the next day's row ("2025-11-04") to the
# "Chart" sheet, mirroring the existing Date/Invalid/Valid rows.
#
# The new date string is written via a temporary formula + "paste values"
# round-trip so Excel's automatic "this looks like a date" literal-input
# heuristic (which would otherwise silently convert a typed
# "2025-11-04" into a date serial number) never fires - the cell ends up
# holding the literal text "2025-11-04", matching the rest of column A.

$wb = $excel.ActiveWorkbook
$chart = $wb.Worksheets.Item("Chart")

$newRow = 31

$chart.Range("A" + $newRow).Formula = '="2025-11-04"'
$chart.Range("A" + $newRow).Copy() | Out-Null
$chart.Range("A" + $newRow).PasteSpecial(-4163) | Out-Null  # xlPasteValues

$chart.Cells.Item($newRow, 2).Value = 0
$chart.Cells.Item($newRow, 3).Value = 105
